$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New homework rows for 2020-02-27 and 2020-02-28
$newRows = @(
    @{ Row = 30; Timestamp = 1582761600; Date = "2020-02-27"; Id = "03022"; Name = "MCOM"; Open = 0.405; High = 0.405; Low = 0.405; Close = 0.405; Vol = "-" },
    @{ Row = 31; Timestamp = 1582848000; Date = "2020-02-28"; Id = "03022"; Name = "MCOM"; Open = 0.405; High = 0.405; Low = 0.405; Close = 0.405; Vol = "-" }
)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Timestamp

    $ws.Cells.Item($r.Row, 2).Value = "'" + $r.Date
    $ws.Cells.Item($r.Row, 2).Style = "Normal"

    $ws.Cells.Item($r.Row, 3).Value = "'" + $r.Id
    $ws.Cells.Item($r.Row, 3).Style = "Normal"

    $ws.Cells.Item($r.Row, 4).Value = $r.Name

    $ws.Cells.Item($r.Row, 5).Value = $r.Open
    $ws.Cells.Item($r.Row, 6).Value = $r.High
    $ws.Cells.Item($r.Row, 7).Value = $r.Low
    $ws.Cells.Item($r.Row, 8).Value = $r.Close

    $ws.Cells.Item($r.Row, 9).Value = "'" + $r.Vol
    $ws.Cells.Item($r.Row, 9).Style = "Normal"
}
